$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original row 18 (week of 2021-11-08) becomes historical data and is
# preserved as new row 19, while row 18 is updated with the latest weekly
# figures (week of 2022-04-05).

# Build new row 19 with the values that row 18 currently holds.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(19, $col).Value2 = $ws.Cells.Item(18, $col).Value2()
}
$ws.Cells.Item(19, 4).NumberFormat = $ws.Cells.Item(18, 4).NumberFormat

# Now update row 18 with this week's figures.
$ws.Cells.Item(18, 4).Value2 = 44656   # D18 Fecha
$ws.Cells.Item(18, 10).Value2 = 85     # J18 Volumen
$ws.Cells.Item(18, 11).Value2 = 5000   # K18 Precio minimo
$ws.Cells.Item(18, 12).Value2 = 5000   # L18 Precio maximo
$ws.Cells.Item(18, 13).Value2 = 5000   # M18 Precio promedio ponderado
$ws.Cells.Item(18, 16).Value2 = 5000   # P18 Precio $/Kg
